# "1. Informacion Demografica.xlsx" - commit "Ya estan todas las tablas"
#
# - Rename the column headers in row 1:
#     D1: "Densidad de Poblacion (hab/km2)" -> "Densidad de poblacion (hab./km2)"
#     E1: "Poblacion Urbana" -> "Poblacion Urbana%"
#     F1: "Poblacion Rural"  -> "Poblacion Rural%"
# - Refresh "Superficie (km2)" (column C) with updated/more precise values for
#   every municipio row (2-85), and recompute the dependent
#   "Densidad de poblacion" (column D = Poblacion total / Superficie) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Densidad de población (hab./km2)"
$ws.Range("E1").Value = "Población Urbana%"
$ws.Range("F1").Value = "Población Rural%"

$ws.Range("C2").Value = 241.60395471
$ws.Range("D2").Value = 92.1673654999916
$ws.Range("C3").Value = 238.8680859
$ws.Range("D3").Value = 192.847026116686
$ws.Range("C4").Value = 271.80287393
$ws.Range("D4").Value = 224.434712988761
$ws.Range("C5").Value = 120.01906514
$ws.Range("D5").Value = 85.9280147530734
$ws.Range("C6").Value = 238.99099886
$ws.Range("D6").Value = 78.9653170622344
$ws.Range("C7").Value = 433.52336674
$ws.Range("D7").Value = 44.2006163222389
$ws.Range("C8").Value = 272.31966558
$ws.Range("D8").Value = 46.070855636808
$ws.Range("C9").Value = 322.22397435
$ws.Range("D9").Value = 144.871281207943
$ws.Range("C10").Value = 137.59729501
$ws.Range("D10").Value = 144.15981068929
$ws.Range("C11").Value = 62.45009249
$ws.Range("D11").Value = 504.803095448546
$ws.Range("C12").Value = 142.64968673
$ws.Range("D12").Value = 138.885688809812
$ws.Range("C13").Value = 458.51095916
$ws.Range("D13").Value = 65.7236199004007
$ws.Range("C14").Value = 122.32060282
$ws.Range("D14").Value = 510.707097249408
$ws.Range("C15").Value = 211.0101449
$ws.Range("D15").Value = 76.5366044729919
$ws.Range("C16").Value = 593.62836729
$ws.Range("D16").Value = 32.7326001766144
$ws.Range("C17").Value = 391.39363559
$ws.Range("D17").Value = 154.373997188072
$ws.Range("C18").Value = 278.33061647
$ws.Range("D18").Value = 46.5884787108847
$ws.Range("C19").Value = 231.63244062
$ws.Range("D19").Value = 98.8764783494772
$ws.Range("C20").Value = 222.82416015
$ws.Range("D20").Value = 84.8606362401227
$ws.Range("C21").Value = 239.53867876
$ws.Range("D21").Value = 10.8249741270302
$ws.Range("C22").Value = 123.01188621
$ws.Range("D22").Value = 123.362062541615
$ws.Range("C23").Value = 142.33870812
$ws.Range("D23").Value = 114.410199552119
$ws.Range("C24").Value = 98.03504701
$ws.Range("D24").Value = 369.745321755214
$ws.Range("C25").Value = 302.84595862
$ws.Range("D25").Value = 58.138467755129
$ws.Range("C26").Value = 292.31124081
$ws.Range("D26").Value = 70.7225625080812
$ws.Range("C27").Value = 107.45597616
$ws.Range("D27").Value = 118.802140711017
$ws.Range("C28").Value = 213.91803679
$ws.Range("D28").Value = 106.797913550542
$ws.Range("C29").Value = 394.04540566
$ws.Range("D29").Value = 321.742109358312
$ws.Range("C30").Value = 660.7457791
$ws.Range("D30").Value = 71.7749571773239
$ws.Range("C31").Value = 486.59532185
$ws.Range("D31").Value = 202.743420600356
$ws.Range("C32").Value = 440.98587366
$ws.Range("D32").Value = 27.8693734518026
$ws.Range("C33").Value = 38.38517377
$ws.Range("D33").Value = 274.14230460575
$ws.Range("C34").Value = 110.7903024
$ws.Range("D34").Value = 26.1304458719484
$ws.Range("C35").Value = 177.20877914
$ws.Range("D35").Value = 53.4623625645277
$ws.Range("C36").Value = 146.3403203
$ws.Range("D36").Value = 89.3670314045363
$ws.Range("C37").Value = 245.72337691
$ws.Range("D37").Value = 38.4538098036185
$ws.Range("C38").Value = 796.90977735
$ws.Range("D38").Value = 26.3041069337935
$ws.Range("C39").Value = 192.14150463
$ws.Range("D39").Value = 46.2055297063279
$ws.Range("C40").Value = 53.42708858
$ws.Range("D40").Value = 268.103697594371
$ws.Range("C41").Value = 232.78663733
$ws.Range("D41").Value = 42.1802561891923
$ws.Range("C42").Value = 135.86509674
$ws.Range("D42").Value = 347.56535072703
$ws.Range("C43").Value = 198.26577265
$ws.Range("D43").Value = 58.3963628479573
$ws.Range("C44").Value = 249.70425397
$ws.Range("D44").Value = 25.0896806938367
$ws.Range("C45").Value = 341.31490835
$ws.Range("D45").Value = 49.6550241005609
$ws.Range("C46").Value = 79.73465508
$ws.Range("D46").Value = 116.574154496236
$ws.Range("C47").Value = 323.98772958
$ws.Range("D47").Value = 118.806968553713
$ws.Range("C48").Value = 385.39864483
$ws.Range("D48").Value = 12.3197111969461
$ws.Range("C49").Value = 154.03052242
$ws.Range("D49").Value = 2040.70592673122
$ws.Range("C50").Value = 188.07657968
$ws.Range("D50").Value = 99.5498750129121
$ws.Range("C51").Value = 90.97783726
$ws.Range("D51").Value = 259.854495468362
$ws.Range("C52").Value = 112.53583047
$ws.Range("D52").Value = 1801.63952363642
$ws.Range("C53").Value = 297.36447474
$ws.Range("D53").Value = 130.785629433389
$ws.Range("C54").Value = 359.39863319
$ws.Range("D54").Value = 49.2461527827882
$ws.Range("C55").Value = 205.71492716
$ws.Range("D55").Value = 178.868886706413
$ws.Range("C56").Value = 256.22499311
$ws.Range("D56").Value = 71.5347858049553
$ws.Range("C57").Value = 64.29112174
$ws.Range("D57").Value = 615.341573288903
$ws.Range("C58").Value = 420.20742226
$ws.Range("D58").Value = 36.0345848213766
$ws.Range("C59").Value = 240.01201236
$ws.Range("D59").Value = 72.6671962311612
$ws.Range("C60").Value = 525.02869222
$ws.Range("D60").Value = 72.3960434224666
$ws.Range("C61").Value = 176.63193089
$ws.Range("D61").Value = 99.0930683473094
$ws.Range("C62").Value = 242.90861996
$ws.Range("D62").Value = 231.547978862429
$ws.Range("C63").Value = 347.34481951
$ws.Range("D63").Value = 89.9250492466342
$ws.Range("C64").Value = 353.42541861
$ws.Range("D64").Value = 256.195494812206
$ws.Range("C65").Value = 147.84988622
$ws.Range("D65").Value = 73.2499718253757
$ws.Range("C66").Value = 44.93408666
$ws.Range("D66").Value = 261.894719014636
$ws.Range("C67").Value = 90.72783863
$ws.Range("D67").Value = 143.638382626376
$ws.Range("C68").Value = 163.34435588
$ws.Range("D68").Value = 337.532323678841
$ws.Range("C69").Value = 265.68034426
$ws.Range("D69").Value = 53.9746364750514
$ws.Range("C70").Value = 76.80893605
$ws.Range("D70").Value = 2191.17733762698
$ws.Range("C71").Value = 31.52423159
$ws.Range("D71").Value = 604.836312839687
$ws.Range("C72").Value = 531.59848768
$ws.Range("D72").Value = 17.0918469682882
$ws.Range("C73").Value = 82.92694023
$ws.Range("D73").Value = 134.009526568541
$ws.Range("C74").Value = 392.11136348
$ws.Range("D74").Value = 96.2022616871292
$ws.Range("C75").Value = 38.98056924
$ws.Range("D75").Value = 734.365879157695
$ws.Range("C76").Value = 128.80692211
$ws.Range("D76").Value = 165.845124237632
$ws.Range("C77").Value = 336.11244819
$ws.Range("D77").Value = 342.465745079848
$ws.Range("C78").Value = 217.41467837
$ws.Range("D78").Value = 774.414134603492
$ws.Range("C79").Value = 135.41019859
$ws.Range("D79").Value = 134.849517910304
$ws.Range("C80").Value = 176.84300888
$ws.Range("D80").Value = 39.6679520690589
$ws.Range("C81").Value = 154.35402674
$ws.Range("D81").Value = 159.853296484204
$ws.Range("C82").Value = 272.69200265
$ws.Range("D82").Value = 139.919761596279
$ws.Range("C83").Value = 105.25738969
$ws.Range("D83").Value = 203.719663418911
$ws.Range("C84").Value = 319.87990261
$ws.Range("D84").Value = 181.024189164517
$ws.Range("C85").Value = 872.49980798
$ws.Range("D85").Value = 45.7616146557539
